# Update the "想去人数" (want-to-go count) column F values across the four
# sheets to match the freshly generated output (gh-pages regeneration at
# commit 456a3b4). Only column F numeric counters change; everything else
# (including column G) stays the same.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value2  = 39
$ws1.Cells.Item(4, 6).Value2  = 13536
$ws1.Cells.Item(5, 6).Value2  = 790
$ws1.Cells.Item(10, 6).Value2 = 1932
$ws1.Cells.Item(13, 6).Value2 = 24613
$ws1.Cells.Item(15, 6).Value2 = 229
$ws1.Cells.Item(18, 6).Value2 = 389
$ws1.Cells.Item(20, 6).Value2 = 330
$ws1.Cells.Item(21, 6).Value2 = 175
$ws1.Cells.Item(22, 6).Value2 = 155
$ws1.Cells.Item(24, 6).Value2 = 249
$ws1.Cells.Item(25, 6).Value2 = 297
$ws1.Cells.Item(26, 6).Value2 = 26
$ws1.Cells.Item(27, 6).Value2 = 1380
$ws1.Cells.Item(28, 6).Value2 = 92
$ws1.Cells.Item(29, 6).Value2 = 385
$ws1.Cells.Item(30, 6).Value2 = 83

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(3, 6).Value2  = 203
$ws2.Cells.Item(6, 6).Value2  = 31
$ws2.Cells.Item(8, 6).Value2  = 93
$ws2.Cells.Item(9, 6).Value2  = 93
$ws2.Cells.Item(12, 6).Value2 = 4

# Sheet 3: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 6).Value2 = 911
$ws3.Cells.Item(3, 6).Value2 = 4629
$ws3.Cells.Item(4, 6).Value2 = 129

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 6).Value2  = 911
$ws4.Cells.Item(3, 6).Value2  = 39
$ws4.Cells.Item(5, 6).Value2  = 13536
$ws4.Cells.Item(6, 6).Value2  = 790
$ws4.Cells.Item(7, 6).Value2  = 4629
$ws4.Cells.Item(11, 6).Value2 = 1932
$ws4.Cells.Item(13, 6).Value2 = 129
$ws4.Cells.Item(14, 6).Value2 = 24613
$ws4.Cells.Item(17, 6).Value2 = 229
$ws4.Cells.Item(18, 6).Value2 = 203
$ws4.Cells.Item(19, 6).Value2 = 203
$ws4.Cells.Item(24, 6).Value2 = 31
$ws4.Cells.Item(26, 6).Value2 = 93
$ws4.Cells.Item(28, 6).Value2 = 389
$ws4.Cells.Item(31, 6).Value2 = 330
$ws4.Cells.Item(32, 6).Value2 = 175
$ws4.Cells.Item(33, 6).Value2 = 155
$ws4.Cells.Item(35, 6).Value2 = 4
$ws4.Cells.Item(36, 6).Value2 = 249
$ws4.Cells.Item(39, 6).Value2 = 297
$ws4.Cells.Item(40, 6).Value2 = 26
$ws4.Cells.Item(42, 6).Value2 = 1380
$ws4.Cells.Item(43, 6).Value2 = 92
$ws4.Cells.Item(45, 6).Value2 = 385
$ws4.Cells.Item(46, 6).Value2 = 83
